$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "A"
$ws.Range("B8").Value = "BB"
$ws.Range("B9").Value = "AAA"
$ws.Range("B10").Value = "AAA"
$ws.Range("B11").Value = "A"
$ws.Range("B15").Value = "BBB"
$ws.Range("B19").Value = "BB"
$ws.Range("B20").Value = "BB"
$ws.Range("B23").Value = "B"
$ws.Range("B29").Value = "B"
$ws.Range("B30").Value = "A"
$ws.Range("B33").Value = "BBB"
$ws.Range("B36").Value = "BBB"
$ws.Range("B46").Value = "BB"
$ws.Range("B47").Value = "BB"
$ws.Range("B50").Value = "A"
$ws.Range("B52").Value = "B"
$ws.Range("B55").Value = "BB"
$ws.Range("B59").Value = "BB"
$ws.Range("B60").Value = "BB"
$ws.Range("B66").Value = "BB"
$ws.Range("B69").Value = "BB"
$ws.Range("B70").Value = "B"
$ws.Range("B72").Value = "A"
$ws.Range("B73").Value = "A"
$ws.Range("B76").Value = "BBB"
$ws.Range("B83").Value = "AA"
$ws.Range("B86").Value = "BBB"
$ws.Range("B90").Value = "BBB"
$ws.Range("B92").Value = "A"
$ws.Range("B95").Value = "A"
$ws.Range("B99").Value = "A"
$ws.Range("B109").Value = "BBB"
$ws.Range("B114").Value = "A"
$ws.Range("B120").Value = "A"
$ws.Range("B122").Value = "BBB"
$ws.Range("B124").Value = "A"
$ws.Range("B125").Value = "A"
$ws.Range("B126").Value = "A"
$ws.Range("B136").Value = "BBB"
$ws.Range("B141").Value = "B"
$ws.Range("B143").Value = "BBB"
$ws.Range("B144").Value = "A"
$ws.Range("B146").Value = "BBB"
$ws.Range("B151").Value = "B"
$ws.Range("B156").Value = "B"
$ws.Range("B157").Value = "B"
$ws.Range("B159").Value = "AA"
$ws.Range("B162").Value = "A"
$ws.Range("B164").Value = "B"
$ws.Range("B165").Value = "B"
$ws.Range("B168").Value = "BB"
$ws.Range("B169").Value = "BBB"
$ws.Range("B173").Value = "B"
$ws.Range("B176").Value = "A"
$ws.Range("B178").Value = "A"
$ws.Range("B184").Value = "BB"
$ws.Range("B185").Value = "BB"
$ws.Range("B194").Value = "BBB"
$ws.Range("B196").Value = "BB"
$ws.Range("B202").Value = "BB"
$ws.Range("B204").Value = "AAA"
$ws.Range("B206").Value = "B"
$ws.Range("B212").Value = "A"
$ws.Range("B214").Value = "BBB"
$ws.Range("B219").Value = "BBB"
$ws.Range("B228").Value = "B"
$ws.Range("B229").Value = "B"
$ws.Range("B230").Value = "B"
$ws.Range("B233").Value = "AA"
$ws.Range("B235").Value = "BBB"
$ws.Range("B236").Value = "BBB"
$ws.Range("B239").Value = "AAA"
$ws.Range("B245").Value = "BBB"
$ws.Range("B247").Value = "BBB"
$ws.Range("B255").Value = "BB"
$ws.Range("B258").Value = "A"
$ws.Range("B259").Value = "BBB"
$ws.Range("B261").Value = "BBB"
$ws.Range("B267").Value = "BBB"
$ws.Range("B271").Value = "BB"
$ws.Range("B272").Value = "BBB"
$ws.Range("B274").Value = "BB"
$ws.Range("B275").Value = "A"
$ws.Range("B284").Value = "BB"
$ws.Range("B286").Value = "BBB"
$ws.Range("B288").Value = "A"
$ws.Range("B291").Value = "A"
$ws.Range("B292").Value = "A"
$ws.Range("B297").Value = "A"
$ws.Range("B301").Value = "A"
$ws.Range("B303").Value = "A"
$ws.Range("B305").Value = "BBB"
$ws.Range("B306").Value = "BB"
$ws.Range("B307").Value = "A"
$ws.Range("B308").Value = "A"
$ws.Range("B309").Value = "BBB"
$ws.Range("B312").Value = "A"
$ws.Range("B313").Value = "BBB"
$ws.Range("B314").Value = "B"
$ws.Range("B319").Value = "BBB"
$ws.Range("B320").Value = "BBB"
$ws.Range("B323").Value = "B"
$ws.Range("B324").Value = "B"
$ws.Range("B330").Value = "A"
$ws.Range("B331").Value = "A"
$ws.Range("B335").Value = "BB"
$ws.Range("B344").Value = "BBB"
$ws.Range("B345").Value = "BBB"
$ws.Range("B347").Value = "BBB"
$ws.Range("B348").Value = "BB"
$ws.Range("B349").Value = "BB"
$ws.Range("B351").Value = "A"
$ws.Range("B353").Value = "BB"
$ws.Range("B358").Value = "BBB"
$ws.Range("B363").Value = "BB"
$ws.Range("B366").Value = "B"
$ws.Range("B367").Value = "B"
$ws.Range("B368").Value = "BB"
$ws.Range("B374").Value = "BB"
$ws.Range("B380").Value = "B"
$ws.Range("B384").Value = "BB"
$ws.Range("B388").Value = "B"
$ws.Range("B393").Value = "BBB"
$ws.Range("B395").Value = "BBB"
$ws.Range("B396").Value = "B"
$ws.Range("B397").Value = "A"
